# Work for the day: add the new "Process and Reactive Schizophrenia" paper
# row to every tracking sheet, plus a new "Diagnosis by Prognosis" citation
# column on the Citations sheet.

$wb = $excel.ActiveWorkbook

$paperTitle = "Process and Reactive Schizophrenia: Some Conceptions and Issues"

# ---------------------------------------------------------------------------
# Summary sheet: full new row (title, author, paper date, start/end dates,
# days taken)
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("A2:G2").Copy()
$wsSummary.Range("A3:G3").PasteSpecial(-4122)

$wsSummary.Range("A3").Value = 1
$wsSummary.Range("B3").Value = $paperTitle
$wsSummary.Range("C3").Value = "Garmezy N"
$wsSummary.Range("D3").Value = "Fall 1970"
$wsSummary.Range("E3").Value = 44466
$wsSummary.Range("F3").Value = 2958465
$wsSummary.Range("G3").Value = 2913999

$wsSummary.Columns.Item(6).ColumnWidth = 20.7

# ---------------------------------------------------------------------------
# Charactheristics / Techniques / Metrics / Problems sheets: just the new
# paper marker row (index + title), the rest of the row stays blank.
# ---------------------------------------------------------------------------
foreach ($name in @("Charactheristics", "Techniques", "Metrics", "Problems")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A2").Copy()
    $ws.Range("A3").PasteSpecial(-4122)
    $ws.Range("A3").Value = 1
    $ws.Range("B3").Value = $paperTitle
}

# ---------------------------------------------------------------------------
# Citations sheet: new "Diagnosis by Prognosis" column, marked "X" for the
# new paper (the older paper keeps its "X" under "Classifying Schizophrenic").
# ---------------------------------------------------------------------------
$wsCitations = $wb.Worksheets.Item("Citations")

$wsCitations.Range("C1").Copy()
$wsCitations.Range("D1").PasteSpecial(-4122)
$wsCitations.Range("D1").Value = "Diagnosis by Prognosis"
$wsCitations.Columns.Item(4).ColumnWidth = 24.0

$wsCitations.Range("A2").Copy()
$wsCitations.Range("A3").PasteSpecial(-4122)
$wsCitations.Range("A3").Value = 1
$wsCitations.Range("B3").Value = $paperTitle

$wsCitations.Range("C2").Copy()
$wsCitations.Range("D3").PasteSpecial(-4122)
$wsCitations.Range("D3").Value = "X"
